$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.218.50'
$ws.Range('E2').Value = '  +0.60%  '
$ws.Range('D3').Value = '1.605.55'
$ws.Range('E3').Value = '  +0.56%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = "'" + '212.56'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.21%  '
$ws.Range('E6').Value = '  -0.03%  '
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('E8').Value = '  +0.32%  '
$ws.Range('E9').Value = '  -0.42%  '
$ws.Range('E10').Value = '  +1.43%  '
$ws.Range('E11').Value = '  -0.21%  '
$ws.Range('D12').Value = '1.829.73'
$ws.Range('E12').Value = '  +0.61%  '
$ws.Range('D13').Value = '1.603.72'
$ws.Range('E13').Value = '  +0.44%  '
$ws.Range('D14').Value = "'" + '4.01'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.38%  '
$ws.Range('D15').Value = "'" + '0.512'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.62%  '
$ws.Range('D16').Value = '26.219.90'
$ws.Range('E16').Value = '  +0.70%  '
$ws.Range('D17').Value = "'" + '61.83'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.30%  '
$ws.Range('E18').Value = '  +1.06%  '
$ws.Range('E19').Value = '  -0.07%  '
$ws.Range('D20').Value = "'" + '200.74'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.41%  '
$ws.Range('D21').Value = "'" + '4.27'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.96%  '
$ws.Range('D22').Value = "'" + '9.27'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.35%  '
$ws.Range('E23').Value = '  +0.61%  '
$ws.Range('E24').Value = '  +2.31%  '
$ws.Range('D25').Value = "'" + '144.43'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.22%  '
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('E27').Value = '  -2.65%  '
$ws.Range('D28').Value = "'" + '15.18'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.18%  '
$ws.Range('D29').Value = "'" + '6.54'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.70%  '
$ws.Range('E30').Value = '  +4.41%  '
$ws.Range('E31').Value = '  +0.43%  '
$ws.Range('E32').Value = '  +2.40%  '
$ws.Range('E33').Value = '  -0.69%  '
$ws.Range('E34').Value = '  +1.59%  '
$ws.Range('E35').Value = '  +0.81%  '
$ws.Range('D36').Value = '1.162.76'
$ws.Range('E36').Value = '  +4.84%  '
$ws.Range('E37').Value = '  +5.20%  '
$ws.Range('E38').Value = '  -0.06%  '
$ws.Range('E39').Value = '  +0.08%  '
$ws.Range('D40').Value = "'" + '0.786'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.09%  '
$ws.Range('E41').Value = '  +1.30%  '
$ws.Range('D42').Value = "'" + '0.785'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.33%  '
$ws.Range('D43').Value = "'" + '5.29'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').Value = '1.741.18'
$ws.Range('D45').Value = "'" + '91.74'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.58%  '
$ws.Range('E46').Value = '  +2.55%  '
$ws.Range('D47').Value = "'" + '54.11'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.48%  '
$ws.Range('E48').Value = '  +0.14%  '
$ws.Range('E49').Value = '  -0.19%  '
$ws.Range('E50').Value = '  -0.01%  '
$ws.Range('D51').Value = '0.0₇0947'
$ws.Range('E51').Value = '  -4.80%  '
